$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.567.50"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.697.02"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'598.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'160.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "2.695.78"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").Value = "'28.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "3.191.47"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "68.437.55"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "2.703.84"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").Value = "'11.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("D20").Value = "'365.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "'7.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("D22").Value = "'4.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "'74.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'9.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "2.830.32"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "'582.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.75%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.25%  "
$ws.Range("D32").Value = "'8.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").Value = "'1.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.28%  "
$ws.Range("D35").Value = "'0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "'160.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'0.379"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Value = "'5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "'2.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("D44").Value = "'17.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("E45").Value = "  -4.68%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'157.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  +5.72%  "
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("D50").Value = "'0.602"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.23%  "
$ws.Range("E51").Value = "  -0.08%  "
